$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 3.55
$ws.Range("Q2").Value = 1.91
$ws.Range("F3").Value = 3.4
$ws.Range("T3").Value = 1.75
$ws.Range("N4").Value = 3.85
$ws.Range("Q4").Value = 1.91
$ws.Range("G5").Value = 1.85
$ws.Range("H5").Value = 4.4
$ws.Range("F6").Value = 1.88
$ws.Range("G6").Value = 1.9
$ws.Range("H6").Value = 4.2
$ws.Range("Q6").Value = 1.63
$ws.Range("T6").Value = 1.61
$ws.Range("X6").Value = 38
$ws.Range("F7").Value = 3.45
$ws.Range("P7").Value = 2.04
$ws.Range("AB7").Value = 1000
$ws.Range("AG7").Value = 1000
$ws.Range("AN7").Value = 44
$ws.Range("H8").Value = 2.36
$ws.Range("I8").Value = 2.44
$ws.Range("J8").Value = 3.75
$ws.Range("K8").Value = 4
$ws.Range("P8").Value = 2.44
$ws.Range("U8").Value = 2.64
$ws.Range("AC8").Value = 9.4
$ws.Range("AD8").Value = 12.5
$ws.Range("AH8").Value = 15
$ws.Range("J9").Value = 7.4
$ws.Range("P9").Value = 3.25
$ws.Range("Q9").Value = 1.39
$ws.Range("T9").Value = 1.93
$ws.Range("U9").Value = 1.97
$ws.Range("X9").Value = 120
$ws.Range("Y9").Value = 470
$ws.Range("AD9").Value = 180
$ws.Range("AF9").Value = 9.800000000000001
$ws.Range("AH9").Value = 60
$ws.Range("AL9").Value = 70
$ws.Range("AN9").Value = 3.2
$ws.Range("I10").Value = 1.79
$ws.Range("P10").Value = 2.4
$ws.Range("Q10").Value = 1.62
$ws.Range("AA10").Value = 21
$ws.Range("H11").Value = 2.16
$ws.Range("F12").Value = 3.2
$ws.Range("G12").Value = 3.6
$ws.Range("U14").Value = 2.18
$ws.Range("F15").Value = 2.2
$ws.Range("AJ15").Value = 1000
$ws.Range("F16").Value = 1.89
$ws.Range("G16").Value = 2.04
$ws.Range("H16").Value = 4
$ws.Range("I16").Value = 4.7
$ws.Range("J16").Value = 3.55
$ws.Range("P16").Value = 1.97
$ws.Range("X16").Value = 19
$ws.Range("H17").Value = 2.36
$ws.Range("Z17").Value = 17
$ws.Range("AD17").Value = 12
$ws.Range("AN17").Value = 1000
$ws.Range("F18").Value = 1.46
$ws.Range("AB18").Value = 1000
$ws.Range("AK18").Value = 1000
$ws.Range("J19").Value = 3.15
$ws.Range("Q19").Value = 2.16
$ws.Range("G20").Value = 1.48
$ws.Range("I20").Value = 9
$ws.Range("AF20").Value = 11
$ws.Range("H21").Value = 1.72
$ws.Range("I21").Value = 1.77
$ws.Range("J21").Value = 4.3
$ws.Range("K21").Value = 4.6
$ws.Range("P21").Value = 2.3
$ws.Range("Q21").Value = 1.7
$ws.Range("T21").Value = 1.72
$ws.Range("U21").Value = 2.22
$ws.Range("X21").Value = 25
$ws.Range("AA21").Value = 21
$ws.Range("AB21").Value = 24
$ws.Range("AC21").Value = 10
$ws.Range("AD21").Value = 12
$ws.Range("AH21").Value = 25
$ws.Range("AM21").Value = 100
$ws.Range("AO21").Value = 8.6
$ws.Range("F22").Value = 1.22
$ws.Range("J22").Value = 7.4
$ws.Range("T22").Value = 2.28
$ws.Range("Y22").Value = 1000
$ws.Range("F23").Value = 2.52
$ws.Range("G23").Value = 2.56
$ws.Range("H23").Value = 3.15
$ws.Range("I23").Value = 3.25
$ws.Range("U23").Value = 2.12
$ws.Range("F24").Value = 4.5
$ws.Range("G24").Value = 4.7
$ws.Range("H24").Value = 1.76
$ws.Range("I24").Value = 1.78
$ws.Range("J24").Value = 4.4
$ws.Range("R24").Value = 1.63
$ws.Range("S24").Value = 2.44
$ws.Range("T24").Value = 1.61
$ws.Range("U24").Value = 2.44
$ws.Range("AA24").Value = 22
$ws.Range("AH24").Value = 19
$ws.Range("AO24").Value = 7.8
$ws.Range("I25").Value = 20
$ws.Range("P25").Value = 2.3
$ws.Range("Q25").Value = 1.65
$ws.Range("U25").Value = 1.57
$ws.Range("AB25").Value = 8.4
$ws.Range("AN25").Value = 4.6
$ws.Range("F26").Value = 3.4
$ws.Range("G26").Value = 3.45
$ws.Range("H26").Value = 2.32
$ws.Range("I26").Value = 2.36
$ws.Range("J26").Value = 3.5
$ws.Range("K26").Value = 3.6
$ws.Range("N26").Value = 3.75
$ws.Range("P26").Value = 1.97
$ws.Range("AB26").Value = 13.5
$ws.Range("AG26").Value = 14.5
$ws.Range("AH26").Value = 22
$ws.Range("AO26").Value = 23
$ws.Range("J27").Value = 3.95
$ws.Range("F28").Value = 8.800000000000001
$ws.Range("G28").Value = 12
$ws.Range("I28").Value = 1.4
$ws.Range("J28").Value = 5.1
$ws.Range("K28").Value = 6
$ws.Range("AF28").Value = 1000
$ws.Range("H29").Value = 1.17
$ws.Range("J29").Value = 7.6
$ws.Range("K29").Value = 9.199999999999999
$ws.Range("O29").Value = 1.17
$ws.Range("P29").Value = 2.64
$ws.Range("Q29").Value = 1.5
$ws.Range("R29").Value = 1.67
$ws.Range("S29").Value = 2.24
$ws.Range("T29").Value = 2.38
$ws.Range("U29").Value = 1.61
$ws.Range("Y29").Value = 10.5
$ws.Range("AB29").Value = 70
$ws.Range("AE29").Value = 15.5
$ws.Range("AG29").Value = 100
$ws.Range("AO29").Value = 3.65
$ws.Range("M32").Value = 1.04
$ws.Range("Q32").Value = 1.64
$ws.Range("AH33").Value = 17
$ws.Range("AM33").Value = 80
$ws.Range("AO33").Value = 42
$ws.Range("H34").Value = 1.73
$ws.Range("Q34").Value = 1.76
$ws.Range("I35").Value = 4.6
$ws.Range("K35").Value = 3.9
$ws.Range("Q35").Value = 2.02
$ws.Range("AB35").Value = 9.199999999999999
$ws.Range("F36").Value = 1.48
$ws.Range("G36").Value = 1.55
$ws.Range("J36").Value = 4.2
$ws.Range("P36").Value = 1.98
$ws.Range("Q36").Value = 1.83
$ws.Range("T36").Value = 1.97
$ws.Range("U36").Value = 1.84
$ws.Range("X36").Value = 17.5
$ws.Range("AA36").Value = 330
$ws.Range("AC36").Value = 10.5
$ws.Range("AD36").Value = 32
$ws.Range("AF36").Value = 9.199999999999999
$ws.Range("AH36").Value = 1000
$ws.Range("AJ36").Value = 16
$ws.Range("AK36").Value = 17.5
$ws.Range("AM36").Value = 180
$ws.Range("AN36").Value = 8.199999999999999
